$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-03 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-04 Tuesday", 2)

$d.Content.Find.Execute("99×21=2079", $true, $false, $false, $false, $false, $true, 1, $false, "99×20=1980", 2)
$d.Content.Find.Execute("88×91=8008", $true, $false, $false, $false, $false, $true, 1, $false, "24×90=2160", 2)
$d.Content.Find.Execute("33×32=1056", $true, $false, $false, $false, $false, $true, 1, $false, "33×37=1221", 2)
$d.Content.Find.Execute("12×69=828", $true, $false, $false, $false, $false, $true, 1, $false, "33×98=3234", 2)
$d.Content.Find.Execute("91×87=7917", $true, $false, $false, $false, $false, $true, 1, $false, "93×78=7254", 2)

$d.Content.Find.Execute("71×89=6319", $true, $false, $false, $false, $false, $true, 1, $false, "75×22=1650", 2)
$d.Content.Find.Execute("88×41=3608", $true, $false, $false, $false, $false, $true, 1, $false, "55×86=4730", 2)
$d.Content.Find.Execute("82×82=6724", $true, $false, $false, $false, $false, $true, 1, $false, "27×78=2106", 2)
$d.Content.Find.Execute("20×23=460", $true, $false, $false, $false, $false, $true, 1, $false, "56×87=4872", 2)
$d.Content.Find.Execute("43×14=602", $true, $false, $false, $false, $false, $true, 1, $false, "64×92=5888", 2)

$d.Content.Find.Execute("28×81=2268", $true, $false, $false, $false, $false, $true, 1, $false, "49×40=1960", 2)
$d.Content.Find.Execute("49×60=2940", $true, $false, $false, $false, $false, $true, 1, $false, "57×16=912", 2)
$d.Content.Find.Execute("32×82=2624", $true, $false, $false, $false, $false, $true, 1, $false, "97×87=8439", 2)
$d.Content.Find.Execute("67×53=3551", $true, $false, $false, $false, $false, $true, 1, $false, "46×57=2622", 2)
$d.Content.Find.Execute("66×44=2904", $true, $false, $false, $false, $false, $true, 1, $false, "42×19=798", 2)

$d.Content.Find.Execute("17×88=1496", $true, $false, $false, $false, $false, $true, 1, $false, "70×25=1750", 2)
$d.Content.Find.Execute("86×25=2150", $true, $false, $false, $false, $false, $true, 1, $false, "28×86=2408", 2)
$d.Content.Find.Execute("97×72=6984", $true, $false, $false, $false, $false, $true, 1, $false, "95×62=5890", 2)
$d.Content.Find.Execute("28×56=1568", $true, $false, $false, $false, $false, $true, 1, $false, "64×62=3968", 2)
$d.Content.Find.Execute("55×73=4015", $true, $false, $false, $false, $false, $true, 1, $false, "83×28=2324", 2)

$d.Content.Find.Execute("62×17=1054", $true, $false, $false, $false, $false, $true, 1, $false, "24×41=984", 2)
$d.Content.Find.Execute("31×47=1457", $true, $false, $false, $false, $false, $true, 1, $false, "92×60=5520", 2)
$d.Content.Find.Execute("14×19=266", $true, $false, $false, $false, $false, $true, 1, $false, "75×16=1200", 2)
$d.Content.Find.Execute("60×47=2820", $true, $false, $false, $false, $false, $true, 1, $false, "90×63=5670", 2)
$d.Content.Find.Execute("95×54=5130", $true, $false, $false, $false, $false, $true, 1, $false, "20×89=1780", 2)
